$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("N2").Value = 3.4
$ws.Range("S2").Value = 3.75
$ws.Range("T2").Value = 1.84
$ws.Range("X2").Value = 16
$ws.Range("Z2").Value = 14
$ws.Range("AA2").Value = 29
$ws.Range("AB2").Value = 13.5
$ws.Range("AC2").Value = 8.199999999999999
$ws.Range("AD2").Value = 11
$ws.Range("AE2").Value = 25
$ws.Range("AF2").Value = 27
$ws.Range("AG2").Value = 16
$ws.Range("AH2").Value = 19
$ws.Range("AI2").Value = 42
$ws.Range("AJ2").Value = 80
$ws.Range("AL2").Value = 1000
$ws.Range("AO2").Value = 19.5

# Row 4 updates
$ws.Range("N4").Value = 1.01
$ws.Range("O4").Value = 1.08
$ws.Range("P4").Value = 1.24
$ws.Range("Q4").Value = 1.08
$ws.Range("S4").Value = 1.08

# Row 9 updates
$ws.Range("G9").Value = 6.6
$ws.Range("H9").Value = 1.8

# Row 10 updates
$ws.Range("J10").Value = 6.2

# Row 12 updates
$ws.Range("Q12").Value = 1.01
